# Auto-generated edit script: refreshes the cryptocurrency Price (column D)
# and Volume(1h) % change (column E) values on Sheet1 to match the latest
# scrape of coinranking.com, per the commit "Updated cryptos list ... with
# GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look numeric (e.g. "69.281.44", "1.00", "0.0720")
# but must stay as literal text, exactly as in the source sheet, so force
# text storage before writing them - otherwise Excel would parse them as
# numbers and mangle formatting such as trailing zeros / thousands dots.
$priceCells = @("D2","D3","D5","D6","D9","D14","D15","D17","D18","D19","D20","D22","D24","D25","D26","D27","D28","D29","D30","D31","D33","D34","D37","D38","D44","D45","D47","D48","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '69.281.44'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '2.471.38'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '559.04'
$ws.Range("E5").Value = '  -1.91%  '
$ws.Range("D6").Value = '163.52'
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("D9").Value = '2.471.58'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("E10").Value = '  -3.79%  '
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("E12").Value = '  -4.13%  '
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").Value = '2.924.61'
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").Value = '69.166.00'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("E16").Value = '  -2.66%  '
$ws.Range("D17").Value = '23.66'
$ws.Range("E17").Value = '  -2.33%  '
$ws.Range("D18").Value = '2.474.89'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '10.80'
$ws.Range("E19").Value = '  -3.61%  '
$ws.Range("D20").Value = '343.54'
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("E21").Value = '  -4.03%  '
$ws.Range("D22").Value = '3.81'
$ws.Range("E22").Value = '  -2.30%  '
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").Value = '1.92'
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").Value = '1.01'
$ws.Range("E25").Value = '  +0.51%  '
$ws.Range("D26").Value = '67.29'
$ws.Range("E26").Value = '  -2.95%  '
$ws.Range("D27").Value = '3.70'
$ws.Range("E27").Value = '  -2.62%  '
$ws.Range("D28").Value = '2.598.88'
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("D30").Value = '8.19'
$ws.Range("E30").Value = '  -5.03%  '
$ws.Range("D31").Value = '0.0₃0822'
$ws.Range("E31").Value = '  -5.77%  '
$ws.Range("E32").Value = '  -4.99%  '
$ws.Range("D33").Value = '439.71'
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("E35").Value = '  -3.85%  '
$ws.Range("E36").Value = '  -5.41%  '
$ws.Range("D37").Value = '156.73'
$ws.Range("E37").Value = '  +2.34%  '
$ws.Range("D38").Value = '19.08'
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  -3.47%  '
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("E42").Value = '  -3.12%  '
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").Value = '37.46'
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("D45").Value = '1.48'
$ws.Range("E45").Value = '  -6.21%  '
$ws.Range("E46").Value = '  +2.74%  '
$ws.Range("D47").Value = '2.08'
$ws.Range("E47").Value = '  -4.94%  '
$ws.Range("D48").Value = '133.68'
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("D50").Value = '0.0720'
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("D51").Value = '0.485'
$ws.Range("E51").Value = '  -3.89%  '
